# HRV Racefields template: fill in this period's turnover figures.
#
#   C19 (Pari-mutuel Turnover)  -> 88.00
#   C20 (Fixed Odds Turnover)   -> 1913.75
#
# Also normalise the formatting of the spacer row (row 23) under the
# totals line back to the sheet's default/"Normal" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C19").Value = 88
$ws.Range("C20").Value = 1913.75

foreach ($addr in @("B23", "C23", "F23", "G23", "I23")) {
    $ws.Range($addr).Style = "Normal"
}
